$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "id"
$ws.Range("B1").Value = "release_date"
$ws.Range("C1").Value = "tassativita"

# Row 2
$ws.Range("A2").Value = 252778
$ws.Range("B2").Value = 45854.58333333334
$ws.Range("C2").Value = 0

# Row 3
$ws.Range("A3").Value = 251889
$ws.Range("B3").Value = 46022.58333333334
$ws.Range("C3").Value = 0

# Row 4
$ws.Range("A4").Value = 252517
$ws.Range("B4").Value = 45855.58333333334
$ws.Range("C4").Value = 0

# Date number format for release_date column values
$ws.Range("B2:B4").NumberFormat = "yyyy-mm-dd h:mm:ss"

# Column widths (ColumnWidth=29.2 serializes to raw OOXML width="30")
$ws.Range("A1:C1").EntireColumn.ColumnWidth = 29.2
